# Release 2.2.0-ballot prep for StructureDefinition-tddui-status-author.xlsx
#
# Updates the IG publisher "Metadata" sheet (version/date/base-definition)
# and the "Elements" sheet (Extension.value[x] Reference(...) cell now
# carries explicit |2.2.0-ballot version pins), mirroring the regenerated
# FHIR IG spreadsheet export.

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsElements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet ---------------------------------------------------

# Version: 2.1.0 -> 2.2.0-ballot
$wsMetadata.Range("B3").Value = "2.2.0-ballot"

# Date: refreshed publication timestamp
$wsMetadata.Range("B8").Value = "2025-12-19T08:32:44+00:00"

# Base Definition: now pins the FHIR core version used
$wsMetadata.Range("B18").Value = "http://hl7.org/fhir/StructureDefinition/Extension|4.0.1"

# --- Elements sheet -----------------------------------------------------

# The longer reference text widens column K's best-fit width
# (139.5390625 -> 158.5546875 character-units in the original export).
$wsElements.Columns.Item(11).ColumnWidth = 157.65

# Extension.value[x] Type(s) cell (K6): the two referenced profile URLs now
# carry an explicit "|2.2.0-ballot" canonical version pin.
$wsElements.Range("K6").Value = "Reference(https://interop.esante.gouv.fr/ig/fhir/tddui/StructureDefinition/tddui-practitioner|2.2.0-ballot|https://interop.esante.gouv.fr/ig/fhir/tddui/StructureDefinition/tddui-practitioner-role|2.2.0-ballot)`n"

# Re-fit row 6 back to its natural (non-custom) height now that the wider
# column lets the wrapped text render on a single visual line again -
# matches the source export, which carries no explicit row height here.
$wsElements.Rows.Item(6).AutoFit()
